$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'61.354.74"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +7.94%  "
$cell = $ws.Range("D3")
$cell.Value = "'3.419.73"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +5.66%  "
$ws.Range("E4").Value = "  +0.05%  "
$cell = $ws.Range("D5")
$cell.Value = "'413.84"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +4.60%  "
$cell = $ws.Range("D6")
$cell.Value = "'123.65"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +14.43%  "
$cell = $ws.Range("D7")
$cell.Value = "'3.414.81"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +5.69%  "
$cell = $ws.Range("D8")
$cell.Value = "'0.579"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("E9").Value = "  +0.00%  "
$cell = $ws.Range("D10")
$cell.Value = "'0.646"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.08%  "
$cell = $ws.Range("D11")
$cell.Value = "'0.113"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +17.61%  "
$cell = $ws.Range("D12")
$cell.Value = "'41.56"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +5.90%  "
$ws.Range("E13").Value = "  -0.63%  "
$cell = $ws.Range("D14")
$cell.Value = "'3.963.78"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'8.45"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.21%  "
$cell = $ws.Range("D16")
$cell.Value = "'19.59"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +3.60%  "
$cell = $ws.Range("D17")
$cell.Value = "'3.386.87"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +4.44%  "
$cell = $ws.Range("D18")
$cell.Value = "'61.406.57"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +8.31%  "
$ws.Range("E19").Value = "  +0.37%  "
$cell = $ws.Range("D20")
$cell.Value = "'10.95"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  +4.65%  "
$cell = $ws.Range("D22")
$cell.Value = "'3.40"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.99%  "
$cell = $ws.Range("D23")
$cell.Value = "'13.05"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$cell = $ws.Range("D24")
$cell.Value = "'298.50"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$cell = $ws.Range("D25")
$cell.Value = "'75.94"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.76%  "
$cell = $ws.Range("D26")
$cell.Value = "'3.13"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "
$cell = $ws.Range("D27")
$cell.Value = "'31.08"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +10.46%  "
$cell = $ws.Range("D28")
$cell.Value = "'8.18"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +13.99%  "
$cell = $ws.Range("D29")
$cell.Value = "'7.73"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "
$cell = $ws.Range("D30")
$cell.Value = "'4.27"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  +0.60%  "
$cell = $ws.Range("D32")
$cell.Value = "'42.95"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +4.36%  "
$ws.Range("E33").Value = "  +5.38%  "
$cell = $ws.Range("D34")
$cell.Value = "'11.45"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +2.51%  "
$cell = $ws.Range("D35")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$cell = $ws.Range("D36")
$cell.Value = "'2.52"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +18.26%  "
$ws.Range("E37").Value = "  -0.81%  "
$cell = $ws.Range("D38")
$cell.Value = "'52.09"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.35%  "
$cell = $ws.Range("D39")
$cell.Value = "'3.57"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.22%  "
$cell = $ws.Range("D40")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$cell = $ws.Range("D41")
$cell.Value = "'3.02"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +2.42%  "
$cell = $ws.Range("D42")
$cell.Value = "'1.98"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +5.19%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Range("D43")
$cell.Value = "'17.67"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.60%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D44")
$cell.Value = "'134.96"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D45")
$cell.Value = "'0.122"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$cell = $ws.Range("D46")
$cell.Value = "'4.00"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.92%  "
$cell = $ws.Range("D47")
$cell.Value = "'0.284"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$cell = $ws.Range("D48")
$cell.Value = "'22.43"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "
$cell = $ws.Range("D49")
$cell.Value = "'2.20"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -3.47%  "
$cell = $ws.Range("D50")
$cell.Value = "'2.202.88"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "
$cell = $ws.Range("D51")
$cell.Value = "'3.761.40"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +5.63%  "
